# fixed #44 TournRPG-44 ダイアログの実装
#
# Adds 3 new rows (5-7) to the "ui" sheet for a confirm/cancel dialog:
#   row5 -> はい
#   row6 -> いいえ
#   row7 -> アイテムを捨てて<val1>を手に入れますか？
# Column A keeps the existing "ROW()-2" running index formula, and row
# styling follows the existing pattern (regular rows use row3/row4's
# style, the new last row inherits the special "last row" style that
# used to sit on row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ui")

# 1) Move the "last row" look (currently A4:B4) down onto the new last
#    row (A7:B7) before anything else overwrites it.
$ws.Range("A4:B4").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)

# 2) Re-stamp row 4 (no longer last) plus the two freshly inserted rows
#    (5 and 6) with the regular "interior row" look taken from row 3.
$ws.Range("A3:B3").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)
$ws.Range("A5:B5").PasteSpecial(-4122)
$ws.Range("A6:B6").PasteSpecial(-4122)

# 3) Match the row height used by the other data rows.
$ws.Rows.Item(5).RowHeight = 20
$ws.Rows.Item(6).RowHeight = 20
$ws.Rows.Item(7).RowHeight = 20

# 4) Extend the running-index formula down through the new rows.
$ws.Range("A5:A7").Formula = "=ROW()-2"

# 5) New dialog strings.
$ws.Range("B5").Value = "はい"
$ws.Range("B6").Value = "いいえ"
$ws.Range("B7").Value = "アイテムを捨てて<val1>を手に入れますか？"
